$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 (shifts rows 3-11 down to 4-12), copying formatting from row above.
$ws.Rows("3:3").Insert()

$ws.Range("A3").Value = "9:00 -10:00"

$ws.Range("A3").Select()
